# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - reuse the existing header style (from AC1)
# so the new headers look consistent with the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the team's season record (74-88-0) for every player row (2-47).
$ws.Range("AD2:AD47").Value = 74
$ws.Range("AE2:AE47").Value = 88
$ws.Range("AF2:AF47").Value = 0
